# Swap the "<X> oddball - <Y> attended" strings in Sheet1!E2:E10 to the
# reworded "<Y> attended - <X> oddball" form, and move the active
# selection from E13 to E10 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New text is written in this specific row order (not top-to-bottom) so
# that new shared-string-table entries are appended in the same sequence
# the original author produced: the three "matched" pairs first, then the
# six "mismatched" pairs.
$ws.Range("E2").Value = "Vibr attended - Vibr oddball"
$ws.Range("E6").Value = "Harm attended - Harm oddball"
$ws.Range("E10").Value = "Keyb attended - Keyb oddball"
$ws.Range("E3").Value = "Vibr attended - Harm oddball"
$ws.Range("E4").Value = "Vibr attended - Keyb oddball"
$ws.Range("E5").Value = "Harm attended - Vibr oddball"
$ws.Range("E7").Value = "Harm attended - Keyb oddball"
$ws.Range("E8").Value = "Keyb attended - Vibr oddball"
$ws.Range("E9").Value = "Keyb attended - Harm oddball"

$ws.Range("E10").Select()
